$d = $word.ActiveDocument

# The JavaScript Promise exercise block ("====", "16. What is the output
# of the following code?", its code listing, "Answer:", "success",
# "error") that used to sit right after the "Success: Hattori" paragraph
# is being removed in its entirety. The blank paragraph that follows the
# exercise is left untouched.

$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("Success: Hattori", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Grow the found hit out to the whole paragraph so we know exactly
    # where the paragraph mark that ends "Success: Hattori" is.
    [void]$anchor.Expand(4)   # wdParagraph
    $deleteStart = $anchor.End
    $deleteEnd = $deleteStart

    # Walk forward paragraph by paragraph until we hit the lone "error"
    # paragraph that ends the exercise (the block's last line), collecting
    # the end boundary as we go.
    $cur = $anchor.Next(4, 1)
    $guard = 0
    while ($cur -ne $null -and $guard -lt 100) {
        $guard = $guard + 1
        $deleteEnd = $cur.End
        if ($cur.Text.Trim() -eq "error") {
            break
        }
        $cur = $cur.Next(4, 1)
    }

    $d.Range($deleteStart, $deleteEnd).Delete()
}
